# feat: add 2022-Q1 data
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert a new "2022-Q1" sheet right before "总计" by duplicating the
#    most recent quarter sheet ("2021-Q4") so all styling/layout carries
#    over, then updating the data that changed for the new quarter.
# ---------------------------------------------------------------------
$src = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")
$src.Copy($totalSheet)
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

$ws = $wb.Worksheets.Item("2022-Q1")

# B2 (fund code "590003") and C2 (fund name) carry over unchanged from
# 2021-Q4; only the numeric-looking metrics & the rank actually change.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "12.43"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "78.81"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "4.95"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "0.6153"
$ws.Range("H2").Value = 7

# Re-apply the (unstyled) General format those cells originally had --
# setting NumberFormat above pulled them into a distinct text style.
$ws.Range("B2").Copy()
$ws.Range("D2:G2").PasteSpecial(-4122)

# Restore the original active tab (copying a sheet makes the copy active).
$wb.Worksheets.Item("2021-Q3").Activate()

# ---------------------------------------------------------------------
# 2) Prepend a "2022-Q1" row to the "总计" summary sheet, pushing the
#    existing rows down and renumbering the index column (A).
# ---------------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")
$tot.Rows.Item(2).Insert()
$tot.Range("A2:D2").ClearFormats()

$tot.Range("A3").Copy()
$tot.Range("A2").PasteSpecial(-4122)

$tot.Range("A2").Value = 0
$tot.Range("B2").Value = "2022-Q1"
$tot.Range("C2").Value = 1
$tot.Range("D2").Value = 0.62

$tot.Range("A3").Value = 1
$tot.Range("A4").Value = 2
